# Weekly fruit/vegetable price update.
# A new weekly record for "Coliflor" (Macroferia Regional de Talca) needs to
# be inserted as row 318, pushing the existing rows 318:342 down to 319:343.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 318 (shifts rows 318:342 -> 319:343,
# inherits formatting - incl. the date number format - from the row above).
$ws.Rows("318:318").Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A318").Value = 5
$ws.Range("B318").Value = "Macroferia Regional de Talca"
$ws.Range("C318").Value = "Maule"
$ws.Range("D318").Value = 44826
$ws.Range("E318").Value = 7
$ws.Range("F318").Value = 100112008
$ws.Range("G318").Value = "Coliflor"
$ws.Range("H318").Value = "Sin especificar"
$ws.Range("I318").Value = "Primera"
$ws.Range("J318").Value = 3000
$ws.Range("K318").Value = 1300
$ws.Range("L318").Value = 1300
$ws.Range("M318").Value = 1300
$ws.Range("N318").Value = "`$/unidad"
$ws.Range("O318").Value = "Región del Maule"
$ws.Range("P318").Value = 1300
$ws.Range("Q318").Value = 1
$ws.Range("R318").Value = "Hortaliza"
